# The "markdown"/export-list sheet drops the extra "Column1" field and
# renames the Loads property from NumPhases to VoltagesMagAng.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Loads row: Property changes from "NumPhases" to "VoltagesMagAng"
$ws.Range("B3").Value = "VoltagesMagAng"

# Drop the third table column ("Column1") from the table definition...
$table = $ws.ListObjects.Item("Table1")
$table.ListColumns.Item("Column1").Delete()

# ...and remove the now-unused worksheet column so it disappears entirely.
$ws.Range("C1:C4").EntireColumn.Delete()

# Leave the selection where the author left it after editing.
$ws.Range("B9").Select()
